$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.066.36"
$ws.Range("E2").Value = '  +3.72%  '
$ws.Range("D3").Value = "'2.583.38"
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'519.72"
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").Value = "'140.02"
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.35%  '
$ws.Range("D9").Value = "'2.597.86"
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("E12").Value = '  +2.76%  '
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("D14").Value = "'3.044.01"
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").Value = "'58.785.40"
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").Value = "'20.48"
$ws.Range("E16").Value = '  +2.64%  '
$ws.Range("D17").Value = "'2.609.62"
$ws.Range("E17").Value = '  +3.30%  '
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").Value = "'338.48"
$ws.Range("E19").Value = '  +2.52%  '
$ws.Range("E20").Value = '  +2.05%  '
$ws.Range("D21").Value = "'10.18"
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").Value = "'6.51"
$ws.Range("E22").Value = '  +6.26%  '
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = "'66.00"
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = "'0.404"
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = "'7.10"
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = "'0.0₃0724"
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("E31").Value = '  -4.87%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = "'1.56"
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'18.74"
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("D34").Value = "'148.67"
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").Value = "'3.99"
$ws.Range("E35").Value = '  +0.73%  '
$ws.Range("D36").Value = "'1.12"
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").Value = "'36.28"
$ws.Range("E37").Value = '  +2.04%  '
$ws.Range("D38").Value = "'0.835"
$ws.Range("E38").Value = '  +1.87%  '
$ws.Range("E39").Value = '  +2.29%  '
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = "'275.27"
$ws.Range("E43").Value = '  +4.22%  '
$ws.Range("D44").Value = "'10.76"
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("D45").Value = "'0.0951"
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = "'0.587"
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").Value = "'0.0521"
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").Value = "'18.59"
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = "'1.981.37"
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = "'4.59"
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = "'0.0220"
$ws.Range("E51").Value = '  +0.04%  '
